# Automatische test-sync: 2025-06-24 20:06:50
#
# A new incoming mail was logged on the "Logs" sheet (row 16) and the
# "Dashboard" pivot-style summary was re-synced to reflect the updated
# category counts (Productinformatie: 1 -> 2), which also changes the
# sort order of the summary rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Append the new mail-log entry as row 16 on the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A16").Value = "Is product Y nog op voorraad?"
$logs.Range("B16").Value = "mailmind.test@zohomail.eu"
$logs.Range("C16").Value = "Ik wil graag weten of product Y beschikbaar is."
$logs.Range("D16").Value = "Productinformatie"
$logs.Range("E16").Value = "Beste klant,`nDank u voor uw interesse in product Y. Helaas is product Y op dit moment niet beschikbaar in onze webshop. Wij raden u aan om regelmatig onze website te bezoeken of u in te schrijven voor onze nieuwsbrief om op de hoogte te blijven van nieuwe voorraad en aanbiedingen. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[E-mailassistent]"
$logs.Range("F16").Value = "2025-06-24 20:06:09"
$logs.Range("G16").Value = "Ja"

# The multi-line answer in E16 would otherwise trigger an explicit
# custom row height; auto-fit the row back to the sheet default so the
# row matches the others (no explicit ht/customHeight attribute).
$logs.Rows.Item(16).AutoFit()

# Extend the existing conditional-formatting rules (Categorie / Beantwoord
# columns) so they keep covering the newly added row.
$cfD = $logs.Range("D2:D15").FormatConditions
for ($i = 1; $i -le $cfD.Count(); $i++) {
    $cfD.Item($i).ModifyAppliesToRange($logs.Range("D2:D16"))
}

$cfG = $logs.Range("G2:G15").FormatConditions
for ($i = 1; $i -le $cfG.Count(); $i++) {
    $cfG.Item($i).ModifyAppliesToRange($logs.Range("G2:G16"))
}

# ---------------------------------------------------------------------
# 2) Re-sync the "Dashboard" category counts.
#    The new row bumps "Productinformatie" from 1 to 2 occurrences,
#    which re-sorts the (count desc) summary table: "Productinformatie"
#    now ties "IT / Technisch probleem" at 2 and moves above it,
#    pushing "Offerte / Prijsaanvraag" (still at 1) down to row 7.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "IT / Technisch probleem"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Offerte / Prijsaanvraag"
$dash.Range("B7").Value = 1
